$wb = $excel.ActiveWorkbook

$updates = @{
    "Auburn Gresham" = @{ 6 = 142; 7 = 537 }
    "Austin" = @{ 2 = 335; 3 = 368; 6 = 415; 7 = 1225 }
    "Belmont Cragin" = @{ 2 = 117; 7 = 347 }
    "Brighton Park" = @{ 3 = 44; 6 = 58; 7 = 183 }
    "By Neighborhood" = @{ 7 = 537; 8 = 1225; 11 = 347; 15 = 183; 19 = 532; 20 = 417; 21 = 56; 22 = 47; 23 = 189; 29 = 968; 30 = 70; 31 = 196; 33 = 770; 34 = 102; 36 = 238; 42 = 665; 43 = 160; 47 = 120; 48 = 225; 52 = 475; 53 = 233; 54 = 355; 63 = 52; 64 = 114; 65 = 412; 67 = 690; 73 = 155; 76 = 249; 77 = 128; 78 = 206; 80 = 65; 84 = 134; 85 = 853; 89 = 263; 91 = 195; 94 = 238; 96 = 194; 97 = 145; 101 = 18093 }
    "Chatham" = @{ 6 = 170; 7 = 532 }
    "Chicago Lawn" = @{ 3 = 135; 7 = 417 }
    "Chinatown" = @{ 5 = 2; 7 = 56 }
    "Citywide Totals" = @{ 2 = 5215; 3 = 5374; 4 = 1118; 5 = 387; 6 = 5999; 7 = 18093 }
    "Clearing" = @{ 3 = 15; 6 = 7; 7 = 47 }
    "Douglas" = @{ 2 = 52; 3 = 68; 6 = 51; 7 = 189 }
    "Englewood" = @{ 2 = 277; 3 = 348; 6 = 268; 7 = 968 }
    "Fuller Park" = @{ 3 = 24; 7 = 70 }
    "Gage Park" = @{ 3 = 46; 7 = 196 }
    "Garfield Park" = @{ 2 = 210; 3 = 286; 6 = 222; 7 = 770 }
    "Garfield Ridge" = @{ 2 = 36; 7 = 102 }
    "Grand Boulevard" = @{ 3 = 69; 7 = 238 }
    "Humboldt Park" = @{ 2 = 176; 3 = 205; 4 = 26; 5 = 7; 7 = 665 }
    "Hyde Park" = @{ 4 = 21; 6 = 64; 7 = 160 }
    "Kenwood" = @{ 3 = 35; 6 = 40; 7 = 120 }
    "Lake View" = @{ 3 = 52; 7 = 225 }
    "Little Village" = @{ 3 = 133; 5 = 17; 7 = 475 }
    "Logan Square" = @{ 3 = 59; 7 = 233 }
    "Loop" = @{ 6 = 189; 7 = 355 }
    "Near South Side" = @{ 3 = 34; 7 = 114 }
    "New City" = @{ 3 = 104; 6 = 155; 7 = 412 }
    "North Lawndale" = @{ 3 = 245; 6 = 194; 7 = 690 }
    "Portage Park" = @{ 2 = 49; 7 = 155 }
    "River North" = @{ 2 = 53; 5 = 2; 7 = 249 }
    "Riverdale" = @{ 2 = 56; 7 = 128 }
    "Rogers Park" = @{ 3 = 46; 6 = 74; 7 = 206 }
    "Rush & Division" = @{ 6 = 30; 7 = 65 }
    "South Deering" = @{ 2 = 43; 7 = 134 }
    "South Shore" = @{ 2 = 285; 3 = 285; 6 = 207; 7 = 853 }
    "Uptown" = @{ 6 = 79; 7 = 263 }
    "Washington Park" = @{ 2 = 49; 7 = 195 }
    "West Loop" = @{ 3 = 45; 6 = 102; 7 = 238 }
    "West Ridge" = @{ 6 = 84; 7 = 194 }
    "West Town" = @{ 3 = 27; 6 = 85; 7 = 145 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsMap = $updates[$sheetName]
    foreach ($row in $rowsMap.Keys) {
        $cellRef = "K" + $row
        $ws.Range($cellRef).Value = $rowsMap[$row]
    }
}

Write-Output "Applied all updates"